$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = 0.62920401764292377
$ws.Range("U2").Value = 0.92189969867474297
$ws.Range("T3").Value = 0.66846239899843707
$ws.Range("BM3").Value = 0.66433089129035217
$ws.Range("P4").Value = 0.88038474369017572
$ws.Range("X4").Value = 0.7072741776874174
$ws.Range("AT5").Value = 0.84263805938140435
$ws.Range("AV5").Value = 0.94420048277835011
$ws.Range("BH5").Value = 0.67578702674797997
$ws.Range("G6").Value = 0.8399651978984245
$ws.Range("AM7").Value = 0.89447670143980529
$ws.Range("BF7").Value = 0.77811218659576586
$ws.Range("E8").Value = 0.75712034441793175
$ws.Range("F8").Value = 0.93541569112509593
$ws.Range("G8").Value = 0.83732763869330773
$ws.Range("AT8").Value = 0.72119013291010403
$ws.Range("BN9").Value = 0.82387127900828827
$ws.Range("L10").Value = 0.8491291323643273
$ws.Range("U10").Value = 0.69835220539776333
$ws.Range("AB11").Value = 0.89334989605810344
$ws.Range("AQ11").Value = 0.91494201556146737
$ws.Range("BE11").Value = 0.96462013384370748
$ws.Range("BG11").Value = 0.61900232371853403
$ws.Range("A12").Value = 0.93586604860607547
$ws.Range("M12").Value = 0.7587056282862692
$ws.Range("O13").Value = 0.72658179110132837
$ws.Range("AG13").Value = 0.99472744027264515
$ws.Range("AZ13").Value = 0.75308658198938094
$ws.Range("B14").Value = 0.89338930329309507
$ws.Range("AA14").Value = 0.73566876608642828
$ws.Range("AY14").Value = 0.67321881976143483
$ws.Range("BH14").Value = 0.81974352255461014
$ws.Range("BI14").Value = 0.92880398031489075
$ws.Range("R16").Value = 0.80635777870291703
$ws.Range("BK16").Value = 0.72228839741883843
$ws.Range("BM16").Value = 0.8439721451104597
$ws.Range("B17").Value = 0.63017391884893026
$ws.Range("H17").Value = 0.60108379458087224
$ws.Range("AQ17").Value = 0.94184526015966663
$ws.Range("BD18").Value = 0.81090119991881815
$ws.Range("G19").Value = 0.87460701918335348
$ws.Range("AS19").Value = 0.81986581828283844
$ws.Range("D20").Value = 0.8460766892043009
$ws.Range("V20").Value = 0.77265338566111785
$ws.Range("AM21").Value = 0.97383402323105817
$ws.Range("X22").Value = 0.97570145524401841
$ws.Range("BE22").Value = 0.955023649538707
$ws.Range("S23").Value = 0.62243633195502812
$ws.Range("AJ23").Value = 0.67861793461642761
$ws.Range("BH23").Value = 0.82554656248822234
$ws.Range("M24").Value = 0.91585607623156262
$ws.Range("AC24").Value = 0.57858715465043664
$ws.Range("BN24").Value = 0.97581488812417327
$ws.Range("AA25").Value = 0.85801116173516334
$ws.Range("AD25").Value = 0.95428596048111369
$ws.Range("AW25").Value = 0.99201043064320515
$ws.Range("X26").Value = 0.98328434339871551
$ws.Range("AS26").Value = 0.77532506265133994
$ws.Range("BJ26").Value = 0.81267891473083165
$ws.Range("BM26").Value = 0.93941957965402489
$ws.Range("AW27").Value = 0.96053516390894189
$ws.Range("O28").Value = 0.74971703796792766
$ws.Range("AD28").Value = 0.99951405414108441
$ws.Range("AX28").Value = 0.88062822286873077
$ws.Range("AE29").Value = 0.93935862213400312
$ws.Range("AF30").Value = 0.85379430643772913
$ws.Range("AI31").Value = 0.81204599916156162
$ws.Range("AU31").Value = 0.94508536881594085
$ws.Range("O32").Value = 0.83544093466450464
$ws.Range("AG32").Value = 0.80652249338395876
$ws.Range("L33").Value = 0.93203590723713858
$ws.Range("AE33").Value = 0.93343899416457721
$ws.Range("AI33").Value = 0.98350831801230165
$ws.Range("BJ33").Value = 0.87016507757677697
$ws.Range("BE34").Value = 0.92533484047431391
$ws.Range("BG34").Value = 0.84503103152419157
$ws.Range("I35").Value = 0.88177712541782305
$ws.Range("AK36").Value = 0.95374909310136369
$ws.Range("AM37").Value = 0.79305878923610051
$ws.Range("H38").Value = 0.93116680180552358
$ws.Range("AX38").Value = 0.70395940662731382
$ws.Range("BP40").Value = 0.97612268658195345
$ws.Range("AC41").Value = 0.99965798957394791
$ws.Range("V42").Value = 0.87935697421715997
$ws.Range("AM42").Value = 0.99599500925156259
$ws.Range("B43").Value = 0.74934022959056135
$ws.Range("H43").Value = 0.791972499707429
$ws.Range("AA43").Value = 0.64287911048124946
$ws.Range("BF43").Value = 0.80714471892980799
$ws.Range("AP44").Value = 0.55456196626014398
$ws.Range("AT44").Value = 0.89713355599247446
$ws.Range("AI45").Value = 0.6110238819478343
$ws.Range("AR45").Value = 0.83035992182957197
$ws.Range("AK46").Value = 0.61248958676983833
$ws.Range("AY46").Value = 0.9995560665894756
$ws.Range("AM47").Value = 0.77698398181329531
$ws.Range("K48").Value = 0.75042069491325214
$ws.Range("AE48").Value = 0.75282713218637942
$ws.Range("AJ48").Value = 0.9081286118970332
$ws.Range("BL48").Value = 0.86428493131425443
$ws.Range("AF49").Value = 0.80396455610658224
$ws.Range("AN49").Value = 0.90758226969920841
$ws.Range("AO49").Value = 0.98622795858435452
$ws.Range("O51").Value = 0.919946754913188
$ws.Range("B52").Value = 0.98030228986320833
$ws.Range("G53").Value = 0.93066809021893604
$ws.Range("S53").Value = 0.77958206940927144
$ws.Range("AR53").Value = 0.96593561788749494
$ws.Range("BP53").Value = 0.56301815788050602
$ws.Range("AB54").Value = 0.81125777032787028
$ws.Range("AQ54").Value = 0.5061083327966629
$ws.Range("A55").Value = 0.94412326198865704
$ws.Range("AB55").Value = 0.72598973678137402
$ws.Range("BC56").Value = 0.96546883952538465
$ws.Range("BG56").Value = 0.94201194861776993
$ws.Range("F58").Value = 0.74831613585283918
$ws.Range("AX59").Value = 0.77324962143608889
$ws.Range("I60").Value = 0.88429084893402232
$ws.Range("K61").Value = 0.85250832169182855
$ws.Range("BH61").Value = 0.85381075522827343
$ws.Range("C62").Value = 0.86622687950085675
$ws.Range("BK62").Value = 0.80564891878318079
$ws.Range("AA63").Value = 0.82561658013968597
$ws.Range("BL63").Value = 0.64770946347648173
$ws.Range("BG64").Value = 0.98554854776844758
$ws.Range("BH64").Value = 0.69927466471189326
$ws.Range("AX65").Value = 0.75951170420249881
$ws.Range("BB65").Value = 0.72566127904327915
$ws.Range("AP66").Value = 0.94225710597432899
$ws.Range("BE66").Value = 0.93678892583023843
$ws.Range("AA67").Value = 0.93644662363247322
$ws.Range("BN67").Value = 0.91728010725757714
$ws.Range("O68").Value = 0.93802651035940343
$ws.Range("AH68").Value = 0.89794208327328962
$ws.Range("BD68").Value = 0.74174101332424258
$ws.Range("BO68").Value = 0.59987326688412912

Write-Output "Updated 136 cells"
